$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the row above (row 17) so the inserted row inherits identical
# formatting/styles, then insert it at row 18 (pushes everything below down by one).
$ws.Rows("17:17").Copy()
$ws.Rows("18:18").Insert()

Write-Output "done"
